$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like label (e.g. "05-10-2021") as plain text, not as
# an auto-converted Excel date. A direct $range.Value = "05-10-2021" gets
# reinterpreted by Excel as a date serial number with a date NumberFormat
# (changing the cell's type/style and styles.xml). Instead, compute the
# label as a text formula result, then Copy / PasteSpecial-values it over
# itself: this "flattens" the formula to its already-text value in place,
# without re-triggering date auto-detection and without touching styles.
function Write-TextLabel($ws, [string]$addr, [string]$text) {
    $r = $ws.Range($addr)
    $r.Formula = '="' + $text + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)   # xlPasteValues
}

# New row 126
Write-TextLabel $ws "A126" "05-10-2021"
$ws.Range("B126").Value = 200000
$ws.Range("D126").Value = 0

# New row 127
Write-TextLabel $ws "A127" "06-10-2021"
$ws.Range("B127").Value = 200000
$ws.Range("D127").Value = 0

# New row 128
Write-TextLabel $ws "A128" "07-10-2021"
$ws.Range("B128").Value = 100000
$ws.Range("C128").Value = 125000
$ws.Range("D128").Value = 50000
$ws.Range("E128").Value = 40000
$ws.Range("F128").Value = 10000
$ws.Range("G128").Value = 2.23
